# Add 40 new Pierce County parcel rows to the ParcelList sheet, then
# update the active sheet / selection state to match the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ParcelList")

$county = "Pierce"
$date = "9/23/2023"

$parcels = @(
    319122025,
    319262013,
    319262020,
    518012029,
    518185010,
    520172039,
    2580000072,
    2930000383,
    3460000800,
    3615220620,
    4015200168,
    4885100520,
    5001930210,
    5002450140,
    5002450160,
    5002520560,
    5003370030,
    5017200560,
    5017201210,
    5017860320,
    5017860840,
    5017880530,
    5018040540,
    5018041090,
    5018060620,
    5018061020,
    5018061030,
    5018120020,
    5018160430,
    5545000290,
    5665000023,
    5670400422,
    5670400480,
    5670400560,
    5820000041,
    6762000810,
    6995100741,
    7470031250,
    7697000330,
    7755000853
)

$startRow = 279
for ($i = 0; $i -lt $parcels.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $county
    # Column B's style (inherited from the <col> default) is a Text format;
    # these parcel numbers were entered as plain numbers (General format,
    # no explicit cell style), so reset to Normal before assigning the value.
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $parcels[$i]
    $ws.Cells.Item($row, 3).Value = $date
}

# Update selections/active views to match the post-edit workbook state.
$auctionNotes = $wb.Worksheets.Item("AuctionNotes")
$auctionNotes.Range("C5").Select()

$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 289
$ws.Range("D306").Select()
